$wb = $excel.ActiveWorkbook

# Add the new "Flights" worksheet after the existing "NewsroomPage" sheet.
$newsroom = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($null, $newsroom)
$ws.Name = "Flights"

# Header + data (shared strings reused / added as needed).
$ws.Range("A1").Value = "Assertions"
$ws.Range("A2").Value = "New York"
$ws.Range("A3").Value = "Mexico"
$ws.Range("A4").Value = "South Africa"
$ws.Range("A5").Value = "Istanbul"
$ws.Range("A6").Value = "Japan"
$ws.Range("A7").Value = "Rome"

# Rows 2-16 carry the explicit black font style (new font/cellXf).
$ws.Range("A2:A16").Font.Color = 0

$ws.Columns("A:A").ColumnWidth = 10

$ws.Range("B9").Select() | Out-Null
